$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Previously the last row (22) used the "last row" date style (YYYY-MM-DD).
# Now that a new last row (23) is appended, row 22 reverts to the regular
# date/time style used by all the other interior rows, and row 23 gets the
# "last row" style instead.
$ws.Range("A22").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Append the new day's data (daily update at 8 AM UTC).
$ws.Range("A23").Value = 45607
$ws.Range("A23").NumberFormat = "YYYY-MM-DD"
$ws.Range("B23").Value = 54
$ws.Range("C23").Value = 48
$ws.Range("D23").Value = 54
